# Make sure all templates have an english language definition.
#
# The "Standard" paragraph style (Word's built-in "Normal" style, which every
# other paragraph style in this template is based on) had no explicit
# language set on its run properties. Setting it to English (US) stamps
# <w:lang w:val="en-US"/> onto that style's <w:rPr>, which is what every
# piece of body text in the document inherits by default.

$d = $word.ActiveDocument

$standard = $d.Styles("Standard")
$standard.LanguageID = "en-US"
